# Bugfixed the naive forecaster component module
# The first data row (old row 2) was a stray/incorrect entry; remove it so
# that every remaining row shifts up by one, and refresh the y_1_forecast
# (column E) values with the corrected forecast numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first data row (date 39400 / 2007). This shifts all
# subsequent rows up by one, which also drops the old trailing row 19.
$ws.Rows.Item(2).Delete()

# Corrected y_1_forecast values (column E) for the remaining 17 data rows.
$E = @(
    "-0.08289353495386509",
    "-0.363786394693788",
    "-0.2139598932957232",
    "0.903223459378788",
    "1.31837503023402",
    "1.658305347589661",
    "1.488472133572305",
    "1.644157643645183",
    "1.586470485311331",
    "1.974604558490256",
    "1.93172124148786",
    "1.008270799755984",
    "-1.119700950349478",
    "0.5759895884974942",
    "0.3530477102890783",
    "-0.01286797263981843",
    "-0.1152140120150968"
)

for ($i = 0; $i -lt $E.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = [double]$E[$i]
}
